$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    "D"  = -0.263
    "G"  = 0.007285714285714285
    "H"  = 0.007285714285714285
    "I"  = -0.06928571428571428
    "J"  = -0.06928571428571428
    "K"  = -12.1
    "L"  = -0.8642857142857142
    "U"  = 2.64
    "V"  = 0.1619631901840491
    "W"  = -0.2694877505567929
    "X"  = 0.07396738683503214
    "Y"  = -0.343455137391825
    "Z"  = -7.106598984771578
    "AA" = 0.4923857868020307
    "AB" = 0.05843459339165084
    "AC" = 0.4339511934103799
    "AD" = 8.109999999999999
    "AF" = 8.109999999999999
    "AG" = 5.469999999999999
    "AH" = 0.3322408848832445
    "AI" = 0.2001974821031844
    "AJ" = 0.2512632062471291
    "AK" = 0.1444415104304198
    "AN" = -18.43181818181818
    "AO" = -2.48081841432225
    "AP" = -12.43181818181818
    "AQ" = -2.48081841432225
}

foreach ($row in @(2, 3)) {
    foreach ($col in $values.Keys) {
        $ws.Range("$col$row").Value = $values[$col]
    }
}
